# Auto-generated edit script: refreshes cached market-price figures
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables, matching
# a scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 315.31818
$ws.Range("I33").Value = 316.75
$ws.Range("K33").Value = 316.75
$ws.Range("M33").Value = -87.75
$ws.Range("H64").Value = 4245.05
$ws.Range("I64").Value = 3428.1428
$ws.Range("K64").Value = 3428.1428
$ws.Range("M64").Value = -3180.1428
$ws.Range("H67").Value = 4245.05
$ws.Range("I67").Value = 3428.1428
$ws.Range("K67").Value = 3428.1428
$ws.Range("M67").Value = -2570.1428
$ws.Range("H112").Value = 2300.9756
$ws.Range("J112").Value = 2449.7026
$ws.Range("L112").Value = 7349.1078
$ws.Range("N112").Value = -9565.1078
$ws.Range("H129").Value = 970.7083
$ws.Range("J129").Value = 1164.8529
$ws.Range("L129").Value = 3494.5587
$ws.Range("N129").Value = -13494.5587
$ws.Range("H134").Value = 117770
$ws.Range("J134").Value = 117770
$ws.Range("L134").Value = 117770
$ws.Range("N134").Value = -127910
$ws.Range("H138").Value = 2473459.5
$ws.Range("I138").Value = 5716836.5
$ws.Range("J138").Value = 5672.674
$ws.Range("K138").Value = 17150509.5
$ws.Range("L138").Value = 17018.022
$ws.Range("M138").Value = -17145369.5
$ws.Range("N138").Value = -27298.022
$ws.Range("H139").Value = 74217.14
$ws.Range("J139").Value = 74217.14
$ws.Range("L139").Value = 74217.14
$ws.Range("N139").Value = -84497.14
$ws.Range("H140").Value = 106151.25
$ws.Range("J140").Value = 107951.43
$ws.Range("L140").Value = 107951.43
$ws.Range("N140").Value = -118311.43

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 50000
$ws.Range("I3").Value = 50000
$ws.Range("K3").Value = 50000
$ws.Range("M3").Value = -49885
$ws.Range("H61").Value = 1850.125
$ws.Range("I61").Value = 1714.3823
$ws.Range("J61").Value = 2619.3333
$ws.Range("K61").Value = 1714.3823
$ws.Range("L61").Value = 2619.3333
$ws.Range("M61").Value = -1502.3823
$ws.Range("N61").Value = -3043.3333
$ws.Range("H63").Value = 10200
$ws.Range("I63").Value = 13500
$ws.Range("J63").Value = 3600
$ws.Range("K63").Value = 13500
$ws.Range("L63").Value = 3600
$ws.Range("M63").Value = -12814
$ws.Range("N63").Value = -4972
$ws.Range("H66").Value = 10200
$ws.Range("I66").Value = 13500
$ws.Range("J66").Value = 3600
$ws.Range("K66").Value = 67500
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = -64068
$ws.Range("N66").Value = -24864
$ws.Range("H92").Value = 275000
$ws.Range("J92").Value = 275000
$ws.Range("L92").Value = 275000
$ws.Range("N92").Value = -279992
$ws.Range("H136").Value = 1850.125
$ws.Range("I136").Value = 1714.3823
$ws.Range("J136").Value = 2619.3333
$ws.Range("K136").Value = 5143.1469
$ws.Range("L136").Value = 7857.999899999999
$ws.Range("M136").Value = -2593.1469
$ws.Range("N136").Value = -12957.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 17100
$ws.Range("I22").Value = 17100
$ws.Range("K22").Value = 17100
$ws.Range("M22").Value = -16927
$ws.Range("H105").Value = 3625.2222
$ws.Range("I105").Value = 3265.875
$ws.Range("K105").Value = 3265.875
$ws.Range("M105").Value = -1518.875
$ws.Range("H132").Value = 76501.336
$ws.Range("J132").Value = 76501.336
$ws.Range("L132").Value = 76501.336
$ws.Range("N132").Value = -86621.336
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 51038
$ws.Range("J138").Value = 51038
$ws.Range("L138").Value = 51038
$ws.Range("N138").Value = -61318

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1736
$ws.Range("I2").Value = 1736
$ws.Range("K2").Value = 1736
$ws.Range("M2").Value = -1623
$ws.Range("H7").Value = 120.833336
$ws.Range("I7").Value = 60
$ws.Range("J7").Value = 181.66667
$ws.Range("K7").Value = 60
$ws.Range("L7").Value = 181.66667
$ws.Range("M7").Value = 53
$ws.Range("N7").Value = -407.66667
$ws.Range("H62").Value = 252002.5
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 252002.5
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14087478
$ws.Range("J131").Value = 16130602
$ws.Range("L131").Value = 48391806
$ws.Range("N131").Value = -48401886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.47059
$ws.Range("I2").Value = 49.5
$ws.Range("J2").Value = 88.28570999999999
$ws.Range("K2").Value = 49.5
$ws.Range("L2").Value = 88.28570999999999
$ws.Range("M2").Value = 63.5
$ws.Range("N2").Value = -314.28571
$ws.Range("H80").Value = 3375.625
$ws.Range("J80").Value = 4750
$ws.Range("L80").Value = 4750
$ws.Range("N80").Value = -6746
$ws.Range("H83").Value = 3375.625
$ws.Range("J83").Value = 4750
$ws.Range("L83").Value = 23750
$ws.Range("N83").Value = -33734
$ws.Range("H140").Value = 50687.5
$ws.Range("J140").Value = 50687.5
$ws.Range("L140").Value = 50687.5
$ws.Range("N140").Value = -61047.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -19988
$ws.Range("H125").Value = 70715
$ws.Range("J125").Value = 70715
$ws.Range("L125").Value = 70715
$ws.Range("N125").Value = -80555
$ws.Range("H132").Value = 3524.842
$ws.Range("I132").Value = 3080.5862
$ws.Range("J132").Value = 4956.3335
$ws.Range("K132").Value = 9241.758600000001
$ws.Range("L132").Value = 14869.0005
$ws.Range("M132").Value = -6711.758600000001
$ws.Range("N132").Value = -19929.0005
$ws.Range("H136").Value = 3587.4656
$ws.Range("I136").Value = 3636.102
$ws.Range("J136").Value = 3322.6667
$ws.Range("K136").Value = 10908.306
$ws.Range("L136").Value = 9968.000100000001
$ws.Range("M136").Value = -8358.306
$ws.Range("N136").Value = -15068.0001
$ws.Range("H139").Value = 79900
$ws.Range("J139").Value = 79900
$ws.Range("L139").Value = 79900
$ws.Range("N139").Value = -90180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 719.875
$ws.Range("I113").Value = 826
$ws.Range("J113").Value = 401.5
$ws.Range("K113").Value = 2478
$ws.Range("L113").Value = 1204.5
$ws.Range("M113").Value = -308
$ws.Range("N113").Value = -5544.5
$ws.Range("H139").Value = 60612.5
$ws.Range("J139").Value = 60612.5
$ws.Range("L139").Value = 60612.5
$ws.Range("N139").Value = -70892.5
